$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1:E1").Font.Bold = $true
$ws.Range("E1").Value = "Size of Original Data File = "
$ws.Range("H1").Value = 286549
$ws.Range("H1").NumberFormat = "#,##0"
